# chore: update Sheets via scheduled runner
# Refreshes cached market-board derived values (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 6951005.5
$ws.Cells.Item(62, 9).Value = 8779113
$ws.Cells.Item(62, 11).Value = 8779113
$ws.Cells.Item(62, 13).Value = -8778489
$ws.Cells.Item(65, 8).Value = 6951005.5
$ws.Cells.Item(65, 9).Value = 8779113
$ws.Cells.Item(65, 11).Value = 43895565
$ws.Cells.Item(65, 13).Value = -43892445
$ws.Cells.Item(100, 8).Value = 2441.36
$ws.Cells.Item(100, 9).Value = 2337.7368
$ws.Cells.Item(100, 10).Value = 2769.5
$ws.Cells.Item(100, 11).Value = 2337.7368
$ws.Cells.Item(100, 12).Value = 2769.5
$ws.Cells.Item(100, 13).Value = -1796.7368
$ws.Cells.Item(100, 14).Value = -3851.5
$ws.Cells.Item(132, 8).Value = 2450.1667
$ws.Cells.Item(132, 9).Value = 1958.8334
$ws.Cells.Item(132, 11).Value = 5876.5002
$ws.Cells.Item(132, 13).Value = -3346.5002
$ws.Cells.Item(137, 8).Value = 46093.2
$ws.Cells.Item(137, 9).Value = 79640.07000000001
$ws.Cells.Item(137, 10).Value = 3397.182
$ws.Cells.Item(137, 11).Value = 238920.21
$ws.Cells.Item(137, 12).Value = 10191.546
$ws.Cells.Item(137, 13).Value = -236370.21
$ws.Cells.Item(137, 14).Value = -15291.546
$ws.Cells.Item(138, 8).Value = 3082.2542
$ws.Cells.Item(138, 10).Value = 3077.0378
$ws.Cells.Item(138, 12).Value = 9231.1134
$ws.Cells.Item(138, 14).Value = -19511.1134

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1654.75
$ws.Cells.Item(32, 9).Value = 1667.8145
$ws.Cells.Item(32, 10).Value = 1232.3334
$ws.Cells.Item(32, 11).Value = 1667.8145
$ws.Cells.Item(32, 12).Value = 1232.3334
$ws.Cells.Item(32, 13).Value = -1380.8145
$ws.Cells.Item(32, 14).Value = -1806.3334
$ws.Cells.Item(61, 8).Value = 2808.6072
$ws.Cells.Item(61, 9).Value = 2060.5652
$ws.Cells.Item(61, 10).Value = 6249.6
$ws.Cells.Item(61, 11).Value = 2060.5652
$ws.Cells.Item(61, 12).Value = 6249.6
$ws.Cells.Item(61, 13).Value = -1848.5652
$ws.Cells.Item(61, 14).Value = -6673.6
$ws.Cells.Item(74, 8).Value = 50558.773
$ws.Cells.Item(74, 9).Value = 58819.793
$ws.Cells.Item(74, 10).Value = 3746.3333
$ws.Cells.Item(74, 11).Value = 58819.793
$ws.Cells.Item(74, 12).Value = 3746.3333
$ws.Cells.Item(74, 13).Value = -57945.793
$ws.Cells.Item(74, 14).Value = -5494.3333
$ws.Cells.Item(77, 8).Value = 50558.773
$ws.Cells.Item(77, 9).Value = 58819.793
$ws.Cells.Item(77, 10).Value = 3746.3333
$ws.Cells.Item(77, 11).Value = 294098.965
$ws.Cells.Item(77, 12).Value = 18731.6665
$ws.Cells.Item(77, 13).Value = -289730.965
$ws.Cells.Item(77, 14).Value = -27467.6665
$ws.Cells.Item(102, 8).Value = 2013.4166
$ws.Cells.Item(102, 9).Value = 1866.2
$ws.Cells.Item(102, 11).Value = 1866.2
$ws.Cells.Item(102, 13).Value = -244.2
$ws.Cells.Item(110, 8).Value = 4122.731
$ws.Cells.Item(110, 9).Value = 4762.722
$ws.Cells.Item(110, 10).Value = 2682.75
$ws.Cells.Item(110, 11).Value = 4762.722
$ws.Cells.Item(110, 12).Value = 2682.75
$ws.Cells.Item(110, 13).Value = -2717.722
$ws.Cells.Item(110, 14).Value = -6772.75
$ws.Cells.Item(136, 8).Value = 2808.6072
$ws.Cells.Item(136, 9).Value = 2060.5652
$ws.Cells.Item(136, 10).Value = 6249.6
$ws.Cells.Item(136, 11).Value = 6181.6956
$ws.Cells.Item(136, 12).Value = 18748.8
$ws.Cells.Item(136, 13).Value = -3631.6956
$ws.Cells.Item(136, 14).Value = -23848.8

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 898.8
$ws.Cells.Item(80, 9).Value = 1231
$ws.Cells.Item(80, 11).Value = 1231
$ws.Cells.Item(80, 13).Value = -233
$ws.Cells.Item(83, 8).Value = 898.8
$ws.Cells.Item(83, 9).Value = 1231
$ws.Cells.Item(83, 11).Value = 6155
$ws.Cells.Item(83, 13).Value = -1163
$ws.Cells.Item(105, 8).Value = 2347.1428
$ws.Cells.Item(105, 9).Value = 2995
$ws.Cells.Item(105, 10).Value = 2088
$ws.Cells.Item(105, 11).Value = 2995
$ws.Cells.Item(105, 12).Value = 2088
$ws.Cells.Item(105, 13).Value = -1248
$ws.Cells.Item(105, 14).Value = -5582
$ws.Cells.Item(107, 8).Value = 621
$ws.Cells.Item(107, 9).Value = 598.75
$ws.Cells.Item(107, 11).Value = 598.75
$ws.Cells.Item(107, 13).Value = 1321.25
$ws.Cells.Item(138, 8).Value = 66730.62
$ws.Cells.Item(138, 10).Value = 66730.62
$ws.Cells.Item(138, 12).Value = 66730.62
$ws.Cells.Item(138, 14).Value = -77010.62

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 271961.9
$ws.Cells.Item(31, 10).Value = 1507.3334
$ws.Cells.Item(31, 12).Value = 1507.3334
$ws.Cells.Item(31, 14).Value = -2097.3334
$ws.Cells.Item(34, 8).Value = 271961.9
$ws.Cells.Item(34, 10).Value = 1507.3334
$ws.Cells.Item(34, 12).Value = 1507.3334
$ws.Cells.Item(34, 14).Value = -1911.3334
$ws.Cells.Item(36, 8).Value = 10000
$ws.Cells.Item(36, 10).Value = 10000
$ws.Cells.Item(36, 12).Value = 10000
$ws.Cells.Item(36, 14).Value = -10776
$ws.Cells.Item(40, 8).Value = 10000
$ws.Cells.Item(40, 10).Value = 10000
$ws.Cells.Item(40, 12).Value = 10000
$ws.Cells.Item(40, 14).Value = -10320
$ws.Cells.Item(64, 8).Value = 33346.08
$ws.Cells.Item(64, 10).Value = 33346.08
$ws.Cells.Item(64, 12).Value = 33346.08
$ws.Cells.Item(64, 14).Value = -33842.08
$ws.Cells.Item(67, 8).Value = 33346.08
$ws.Cells.Item(67, 10).Value = 33346.08
$ws.Cells.Item(67, 12).Value = 33346.08
$ws.Cells.Item(67, 14).Value = -35062.08
$ws.Cells.Item(68, 8).Value = 31583
$ws.Cells.Item(68, 10).Value = 31583
$ws.Cells.Item(68, 12).Value = 31583
$ws.Cells.Item(68, 14).Value = -33081
$ws.Cells.Item(71, 8).Value = 31583
$ws.Cells.Item(71, 10).Value = 31583
$ws.Cells.Item(71, 12).Value = 94749
$ws.Cells.Item(71, 14).Value = -102237
$ws.Cells.Item(132, 8).Value = 4476.361
$ws.Cells.Item(132, 9).Value = 4342.1875
$ws.Cells.Item(132, 11).Value = 13026.5625
$ws.Cells.Item(132, 13).Value = -10496.5625
$ws.Cells.Item(134, 8).Value = 3015.8235
$ws.Cells.Item(134, 9).Value = 2891.8125
$ws.Cells.Item(134, 10).Value = 5000
$ws.Cells.Item(134, 11).Value = 8675.4375
$ws.Cells.Item(134, 12).Value = 15000
$ws.Cells.Item(134, 13).Value = -6140.4375
$ws.Cells.Item(134, 14).Value = -20070
$ws.Cells.Item(139, 8).Value = 70000
$ws.Cells.Item(139, 10).Value = 70000
$ws.Cells.Item(139, 12).Value = 70000
$ws.Cells.Item(139, 14).Value = -80280

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 8040.385
$ws.Cells.Item(55, 9).Value = 1775
$ws.Cells.Item(55, 10).Value = 9179.546
$ws.Cells.Item(55, 11).Value = 5325
$ws.Cells.Item(55, 12).Value = 27538.638
$ws.Cells.Item(55, 13).Value = -5148
$ws.Cells.Item(55, 14).Value = -27892.638
$ws.Cells.Item(68, 8).Value = 3206903.8
$ws.Cells.Item(68, 10).Value = 2185.1052
$ws.Cells.Item(68, 12).Value = 6555.3156
$ws.Cells.Item(68, 14).Value = -8177.3156
$ws.Cells.Item(71, 8).Value = 3206903.8
$ws.Cells.Item(71, 10).Value = 2185.1052
$ws.Cells.Item(71, 12).Value = 19665.9468
$ws.Cells.Item(71, 14).Value = -27777.9468
$ws.Cells.Item(76, 8).Value = 125227500
$ws.Cells.Item(76, 9).Value = 250450000
$ws.Cells.Item(76, 10).Value = 5000
$ws.Cells.Item(76, 11).Value = 751350000
$ws.Cells.Item(76, 12).Value = 15000
$ws.Cells.Item(76, 13).Value = -751349617
$ws.Cells.Item(76, 14).Value = -15766
$ws.Cells.Item(79, 8).Value = 125227500
$ws.Cells.Item(79, 9).Value = 250450000
$ws.Cells.Item(79, 10).Value = 5000
$ws.Cells.Item(79, 11).Value = 751350000
$ws.Cells.Item(79, 12).Value = 15000
$ws.Cells.Item(79, 13).Value = -751348674
$ws.Cells.Item(79, 14).Value = -17652
$ws.Cells.Item(131, 8).Value = 1528.5555
$ws.Cells.Item(131, 10).Value = 1682.7778
$ws.Cells.Item(131, 12).Value = 5048.3334
$ws.Cells.Item(131, 14).Value = -15128.3334

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2559.2
$ws.Cells.Item(80, 10).Value = 3403
$ws.Cells.Item(80, 12).Value = 3403
$ws.Cells.Item(80, 14).Value = -5399
$ws.Cells.Item(83, 8).Value = 2559.2
$ws.Cells.Item(83, 10).Value = 3403
$ws.Cells.Item(83, 12).Value = 17015
$ws.Cells.Item(83, 14).Value = -26999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1789.0555
$ws.Cells.Item(82, 9).Value = 1463.7333
$ws.Cells.Item(82, 11).Value = 1463.7333
$ws.Cells.Item(82, 13).Value = -1102.7333
$ws.Cells.Item(85, 8).Value = 1789.0555
$ws.Cells.Item(85, 9).Value = 1463.7333
$ws.Cells.Item(85, 11).Value = 1463.7333
$ws.Cells.Item(85, 13).Value = -215.7333000000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 314697.8
$ws.Cells.Item(136, 10).Value = 988.5
$ws.Cells.Item(136, 12).Value = 2965.5
$ws.Cells.Item(136, 14).Value = -8065.5
